# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement rows for worker 1047435144 (JAVIER ENRIQUE ALVAREZ CORENA)
# periods 1906 -> 1901 (descending), rows 16-21
$ws.Cells.Item(16, 3).Value = "1047435144"
$ws.Cells.Item(16, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(16, 5).Value = "1906"
$ws.Cells.Item(16, 6).Value = 21534

$ws.Cells.Item(17, 3).Value = "1047435144"
$ws.Cells.Item(17, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(17, 5).Value = "1905"
$ws.Cells.Item(17, 6).Value = 38000

$ws.Cells.Item(18, 3).Value = "1047435144"
$ws.Cells.Item(18, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(18, 5).Value = "1904"
$ws.Cells.Item(18, 6).Value = 38000

$ws.Cells.Item(19, 3).Value = "1047435144"
$ws.Cells.Item(19, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(19, 5).Value = "1903"
$ws.Cells.Item(19, 6).Value = 38000

$ws.Cells.Item(20, 3).Value = "1047435144"
$ws.Cells.Item(20, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(20, 5).Value = "1902"
$ws.Cells.Item(20, 6).Value = 38000

$ws.Cells.Item(21, 3).Value = "1047435144"
$ws.Cells.Item(21, 4).Value = "JAVIER ENRIQUE ALVAREZ CORENA"
$ws.Cells.Item(21, 5).Value = "1901"
$ws.Cells.Item(21, 6).Value = 38000

# New account-statement rows for worker 1044916854 (JUAN CARLOS ARRIETA BABILONIA)
# periods 1906 -> 1901 (descending), rows 22-27
$ws.Cells.Item(22, 3).Value = "1044916854"
$ws.Cells.Item(22, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(22, 5).Value = "1906"
$ws.Cells.Item(22, 6).Value = 21534

$ws.Cells.Item(23, 3).Value = "1044916854"
$ws.Cells.Item(23, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(23, 5).Value = "1905"
$ws.Cells.Item(23, 6).Value = 38000

$ws.Cells.Item(24, 3).Value = "1044916854"
$ws.Cells.Item(24, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(24, 5).Value = "1904"
$ws.Cells.Item(24, 6).Value = 38000

$ws.Cells.Item(25, 3).Value = "1044916854"
$ws.Cells.Item(25, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(25, 5).Value = "1903"
$ws.Cells.Item(25, 6).Value = 38000

$ws.Cells.Item(26, 3).Value = "1044916854"
$ws.Cells.Item(26, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(26, 5).Value = "1902"
$ws.Cells.Item(26, 6).Value = 38000

$ws.Cells.Item(27, 3).Value = "1044916854"
$ws.Cells.Item(27, 4).Value = "JUAN CARLOS ARRIETA BABILONIA"
$ws.Cells.Item(27, 5).Value = "1901"
$ws.Cells.Item(27, 6).Value = 38000

$wb.Save()
